$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 entirely (shrinks used range from A1:K7 to A1:K6)
$ws.Rows.Item(7).Delete()

# Several cells hold digit-only text ("1", "111111111111111111111111", ...)
# that must stay text (matches the original file's t="str" cells) rather
# than be coerced into floating-point numbers by COM's type sniffing.
# Forcing the Text number format before assignment preserves the exact
# digit string (no precision loss, no scientific notation).
$textCells = "C2","D2","E2","C3","D3","E3","B4","C4","E4","C5","D5","E5"
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("A2").Value = "HAYLALA ONE"
$ws.Range("B2").Value = "BG12456"
$ws.Range("C2").Value = "111111111111111111111111"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "11"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "949/DR"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 4334.4
$ws.Range("J2").Value = 433.44
$ws.Range("K2").Value = 3900.96

# Row 3
$ws.Range("A3").Value = "HAYLAL TWO"
$ws.Range("B3").Value = "BG196435"
$ws.Range("C3").Value = "114655555555555555555555"
$ws.Range("D3").Value = "5"
$ws.Range("E3").Value = "5"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "949/DR"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 5665.73
$ws.Range("J3").Value = 566.57
$ws.Range("K3").Value = 5099.16

# Row 4
$ws.Range("A4").Value = "ALI EXPRESSE"
$ws.Range("B4").Value = "11986345"
$ws.Range("C4").Value = "114684354634563543243543"
# D4 is blank text (matches original's empty t="str" cells); a plain ""
# assignment clears the cell entirely, so use the quote-prefix trick to
# commit an empty *text* value instead.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "354354"
$ws.Range("F4").Value = "Logement de fonction"
$ws.Range("G4").Value = "001/LF/TEST DR/AV1"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 30000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30000

# Row 5
$ws.Range("A5").Value = "TETS TESTS"
$ws.Range("B5").Value = "BG432432"
$ws.Range("C5").Value = "321321321312111111111111"
$ws.Range("D5").Value = "11"
$ws.Range("E5").Value = "111"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "001/TEST DR"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 20000
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 17000

# Row 6
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 60000.13
$ws.Range("J6").Value = 4000.01
$ws.Range("K6").Value = 56000.12
